# Applies the "Usiwal" -> "Gráfica Vektra" keyword-report refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title / site rows.
$ws.Range("A1").Value = "Gráfica Vektra"
$ws.Range("A2").Value = "graficavektra.com.br"

# Update existing keyword rows (text + position values).
$ws.Range("A5").Value = "Gráfica para convites especiais"
$ws.Range("C5").Value = 2

$ws.Range("A6").Value = "Gráfica ecológica"
$ws.Range("C6").Value = 4

# Add the two new keyword rows.
$ws.Range("A7").Value = "Gráfica digital"
$ws.Range("B7").Value = "Não encontrado"
$ws.Range("C7").Value = "Não encontrado"

$ws.Range("A8").Value = "Gráfica digital em sp"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 8
